# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp shown in A1
$ws.Range("A1").Value = "Datos actualizados a 17 de Octubre de 2020 a las 12:21"

# Row 20 - Banglades
$ws.Range("B20").Value = 387295
$ws.Range("C20").Value = 1209
$ws.Range("D20").Value = 302298
$ws.Range("E20").Value = 79351
$ws.Range("G20").Value = 23
$ws.Range("H20").Value = 5646

# Row 38 - Catar
$ws.Range("B38").Value = 129227
$ws.Range("C38").Value = 235
$ws.Range("D38").Value = 126218
$ws.Range("E38").Value = 2786
$ws.Range("G38").Value = 1
$ws.Range("H38").Value = 223

# Row 127 - Sri Lanka
$ws.Range("D127").Value = 3395
$ws.Range("E127").Value = 1946

# Row 131 - Hong Kong
$ws.Range("B131").Value = 5238
$ws.Range("C131").Value = 17
$ws.Range("D131").Value = 4963
$ws.Range("E131").Value = 170

# Row 149 - Mali
$ws.Range("B149").Value = 3392
$ws.Range("C149").Value = 188
$ws.Range("D149").Value = 1329
$ws.Range("E149").Value = 2020
$ws.Range("G149").Value = 1
$ws.Range("H149").Value = 43

# Row 150 - Principado de Andorra
$ws.Range("B150").Value = 3378
$ws.Range("D150").Value = 2563
$ws.Range("E150").Value = 683
$ws.Range("H150").Value = 132

# Row 151 - Letonia
$ws.Range("B151").Value = 3377
$ws.Range("D151").Value = 2057
$ws.Range("E151").Value = 1261
$ws.Range("H151").Value = 59
